$wb = $excel.ActiveWorkbook

$wsGeneral = $wb.Worksheets.Item("general")
$wsGeneral.Range("B3").Value = 181.9610138849955
$wsGeneral.Range("B4").Value = 0.01700019836425781
$wsGeneral.Range("B6").Value = 32.48101388499555

$wsX = $wb.Worksheets.Item("x")
$wsX.Range("B2").Value = 3
$wsX.Range("B3").Value = 4
$wsX.Range("B7").Value = 6
$wsX.Range("B10").Value = 5
$wsX.Range("B13").Value = 11

$wsTBar = $wb.Worksheets.Item("TBar")
$wsTBar.Range("B5").Value = 30
$wsTBar.Range("B7").Value = 34.16886835983306
$wsTBar.Range("B8").Value = 30.34885527085025
$wsTBar.Range("B13").Value = 37.87444125446785
$wsTBar.Range("B15").Value = 35.35398438790794

$wsQ = $wb.Worksheets.Item("Q")
$wsQ.Range("C17").Value = 46.91999999999942
$wsQ.Range("C18").Value = 36.10499999999942
$wsQ.Range("C19").Value = 34.91499999999942
$wsQ.Range("C20").Value = 37.48999999999942
$wsQ.Range("C21").Value = 39.43499999999941

$wsQ.Range("C27").Value = 224.1799999999995
$wsQ.Range("C28").Value = 224.6649999999995
$wsQ.Range("C29").Value = 201.1149999999995
$wsQ.Range("C30").Value = 218.9699999999995
$wsQ.Range("C31").Value = 207.1049999999995

$wsQ.Range("C32").Value = 154.3
$wsQ.Range("C33").Value = 148.3449999999993
$wsQ.Range("C34").Value = 128.7049999999993
$wsQ.Range("C35").Value = 146.3249999999992
$wsQ.Range("C36").Value = 134.2149999999993

$wsQ.Range("C43").Value = 159.2149999999989
$wsQ.Range("C44").Value = 142.1399999999989

$wsQ.Range("C48").Value = 247.1799999999994

$wsQ.Range("C57").Value = 226.0399999999994
$wsQ.Range("C58").Value = 247.1799999999994
$wsQ.Range("C59").Value = 221.8549999999994
$wsQ.Range("C60").Value = 238.4549999999994
$wsQ.Range("C61").Value = 224.4749999999994

$wsQ.Range("C67").Value = 224.1799999999995
$wsQ.Range("C68").Value = 224.6649999999995
$wsQ.Range("C69").Value = 201.1149999999995
$wsQ.Range("C70").Value = 218.9699999999995
$wsQ.Range("C71").Value = 207.1049999999995
